$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.598.54'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.415.91'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.38%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '508.82'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.13'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.65%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.993'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.63%  '
$ws.Range('E8').Value = '  -0.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.456.24'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('E11').Value = '  -1.38%  '
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.64'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.851.57'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '57.411.81'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.70%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.96'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.17%  '
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.466.89'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.36'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.14'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '315.72'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.37'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.46%  '
$ws.Range('E23').Value = '  -0.41%  '
$ws.Range('E24').Value = '  -1.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.41'
$ws.Range('D25').Style = 'Normal'
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.559.34'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.995'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.384'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.86%  '
$ws.Range('E29').Value = '  -2.00%  '
$ws.Range('E30').Value = '  +5.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '173.83'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0₃0740'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.51%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.23'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.95%  '
$ws.Range('E35').Value = '  +0.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.997'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.990'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.02'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.11%  '
$ws.Range('E39').Value = '  +5.58%  '
$ws.Range('E40').Value = '  +3.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.826'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.49'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.56%  '
$ws.Range('E43').Value = '  +1.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '134.96'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.86%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.42'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.03'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '260.88'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.573'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.26%  '
$ws.Range('E49').Value = '  -0.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0498'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('E51').Value = '  +1.77%  '
